# The tester filled in the login username field on the "LogIn" sheet with
# "adminYadhu" (replacing the previous "admin" value in B1) and left that
# sheet/cell as the active selection when the workbook was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LogIn")

# Bring the LogIn sheet to the front / make it the active tab (moves
# tabSelected away from whatever sheet had it before, e.g. AdminUsersTest).
$ws.Activate()

# Update the username cell with the new value.
$ws.Range("B1").Value = "adminYadhu"

# Leave the selection on B1, matching the saved workbook state.
$ws.Range("B1").Select()
